$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 data updates (testcase values modified + exam search number added)
$ws.Range("A2").Value = "MXato821"
$ws.Range("B2").Value = 23071960
$ws.Range("C2").Value = "hjuzqjd22"
$ws.Range("D2").Value = "sH#c5%2D"
$ws.Range("F2").Value = "ppxFBVUI"
$ws.Range("G2").Value = "ZlSl"
